{"js": "// Update the date label and every two-digit-by-two-digit multiplication\n// answer in the practice-sheet table. Each old value is unique in the\n// document, so a targeted search-and-replace per pair is safe and order\n// independent.\nconst replacements = [\n  [\"2024-01-21 Sunday\", \"2024-01-22 Monday\"],\n  [\"96\u00d757=5472\", \"94\u00d772=6768\"],\n  [\"38\u00d723=874\", \"51\u00d749=2499\"],\n  [\"26\u00d740=1040\", \"49\u00d742=2058\"],\n  [\"28\u00d732=896\", \"96\u00d760=5760\"],\n  [\"88\u00d715=1320\", \"78\u00d752=4056\"],\n  [\"92\u00d752=4784\", \"11\u00d734=374\"],\n  [\"66\u00d724=1584\", \"25\u00d783=2075\"],\n  [\"85\u00d791=7735\", \"17\u00d766=1122\"],\n  [\"92\u00d779=7268\", \"19\u00d770=1330\"],\n  [\"68\u00d729=1972\", \"63\u00d765=4095\"],\n  [\"65\u00d758=3770\", \"69\u00d743=2967\"],\n  [\"46\u00d733=1518\", \"47\u00d797=4559\"],\n  [\"66\u00d725=1650\", \"32\u00d717=544\"],\n  [\"74\u00d767=4958\", \"32\u00d746=1472\"],\n  [\"64\u00d752=3328\", \"75\u00d781=6075\"],\n  [\"15\u00d734=510\", \"46\u00d714=644\"],\n  [\"28\u00d792=2576\", \"90\u00d755=4950\"],\n  [\"68\u00d754=3672\", \"26\u00d750=1300\"],\n  [\"45\u00d739=1755\", \"47\u00d744=2068\"],\n  [\"89\u00d771=6319\", \"90\u00d735=3150\"],\n  [\"30\u00d745=1350\", \"66\u00d748=3168\"],\n  [\"27\u00d793=2511\", \"19\u00d755=1045\"],\n  [\"45\u00d719=855\", \"24\u00d773=1752\"],\n  [\"34\u00d732=1088\", \"57\u00d758=3306\"],\n  [\"41\u00d731=1271\", \"19\u00d775=1425\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and every two-digit-by-two-digit multiplication\n# answer in the practice-sheet table. Each old value is unique in the\n# document, so Find/Replace (wdReplaceAll) per pair is safe regardless of\n# which paragraph/cell it lives in.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-21 Sunday\", \"2024-01-22 Monday\"),\n    @(\"96\u00d757=5472\", \"94\u00d772=6768\"),\n    @(\"38\u00d723=874\", \"51\u00d749=2499\"),\n    @(\"26\u00d740=1040\", \"49\u00d742=2058\"),\n    @(\"28\u00d732=896\", \"96\u00d760=5760\"),\n    @(\"88\u00d715=1320\", \"78\u00d752=4056\"),\n    @(\"92\u00d752=4784\", \"11\u00d734=374\"),\n    @(\"66\u00d724=1584\", \"25\u00d783=2075\"),\n    @(\"85\u00d791=7735\", \"17\u00d766=1122\"),\n    @(\"92\u00d779=7268\", \"19\u00d770=1330\"),\n    @(\"68\u00d729=1972\", \"63\u00d765=4095\"),\n    @(\"65\u00d758=3770\", \"69\u00d743=2967\"),\n    @(\"46\u00d733=1518\", \"47\u00d797=4559\"),\n    @(\"66\u00d725=1650\", \"32\u00d717=544\"),\n    @(\"74\u00d767=4958\", \"32\u00d746=1472\"),\n    @(\"64\u00d752=3328\", \"75\u00d781=6075\"),\n    @(\"15\u00d734=510\", \"46\u00d714=644\"),\n    @(\"28\u00d792=2576\", \"90\u00d755=4950\"),\n    @(\"68\u00d754=3672\", \"26\u00d750=1300\"),\n    @(\"45\u00d739=1755\", \"47\u00d744=2068\"),\n    @(\"89\u00d771=6319\", \"90\u00d735=3150\"),\n    @(\"30\u00d745=1350\", \"66\u00d748=3168\"),\n    @(\"27\u00d793=2511\", \"19\u00d755=1045\"),\n    @(\"45\u00d719=855\", \"24\u00d773=1752\"),\n    @(\"34\u00d732=1088\", \"57\u00d758=3306\"),\n    @(\"41\u00d731=1271\", \"19\u00d775=1425\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1        # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
